$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.6122786519329679
$ws.Range("C2").Value = 0.1885470363573631
$ws.Range("D2").Value = 0.02063120273205499
$ws.Range("F2").Value = 0.3281835288709161
$ws.Range("G2").Value = 0.1861946332632058
$ws.Range("H2").Value = 0.3542096346093899
$ws.Range("M2").Value = 0.8951897497989307
$ws.Range("O2").Value = 0.9865177010593413
# Row 3
$ws.Range("B3").Value = 0.5345322749281252
$ws.Range("C3").Value = 0.1751280124174173
$ws.Range("D3").Value = 0.01800633603514257
$ws.Range("F3").Value = 0.3256611612239553
$ws.Range("G3").Value = 0.1847294019306247
$ws.Range("H3").Value = 0.3573161180555786
$ws.Range("M3").Value = 0.796492303671414
$ws.Range("O3").Value = 0.9896594350591812
# Row 4
$ws.Range("B4").Value = 0.48657823175364
$ws.Range("C4").Value = 0.1668668200906325
$ws.Range("D4").Value = 0.01638785620799155
$ws.Range("F4").Value = 0.3244310490701565
$ws.Range("G4").Value = 0.1840714663983221
$ws.Range("H4").Value = 0.3594608070543899
$ws.Range("M4").Value = 0.7364158086837307
$ws.Range("O4").Value = 0.9925888887174636
# Row 5
$ws.Range("B5").Value = 0.4669837857179004
$ws.Range("C5").Value = 0.1634952817321675
$ws.Range("D5").Value = 0.0157266590688323
$ws.Range("F5").Value = 0.3240096576535052
$ws.Range("G5").Value = 0.1838638036605857
$ws.Range("H5").Value = 0.360394384645609
$ws.Range("M5").Value = 0.7120598214077205
$ws.Range("O5").Value = 0.9940333706022102
# Row 6
$ws.Range("B6").Value = 0.4637270108456732
$ws.Range("C6").Value = 0.1629351498364713
$ws.Range("D6").Value = 0.01561676959285307
$ws.Range("F6").Value = 0.3239445042005897
$ws.Range("G6").Value = 0.1838329628014961
$ws.Range("H6").Value = 0.360553002085453
$ws.Range("M6").Value = 0.7080229417154982
$ws.Range("O6").Value = 0.9942883425900249
# Row 7
$ws.Range("B7").Value = 0.4863141854037565
$ws.Range("C7").Value = 0.1668213701223351
$ws.Range("D7").Value = 0.01637894569185505
$ws.Range("F7").Value = 0.3244250428805131
$ws.Range("G7").Value = 0.1840684214619372
$ws.Range("H7").Value = 0.3594731563881126
$ws.Range("M7").Value = 0.7360868344652971
$ws.Range("O7").Value = 0.9926073555503336
# Row 8
$ws.Range("B8").Value = 0.5855179442698102
$ws.Range("C8").Value = 0.1839249646226051
$ws.Range("D8").Value = 0.01972759344632635
$ws.Range("F8").Value = 0.3272475424687826
$ws.Range("G8").Value = 0.1856390515843813
$ws.Range("H8").Value = 0.3552314568356962
$ws.Range("M8").Value = 0.8610463363554146
$ws.Range("O8").Value = 0.9873928167925357
# Row 9
$ws.Range("B9").Value = 0.7782568001201753
$ws.Range("C9").Value = 0.2172730665704989
$ws.Range("D9").Value = 0.02623813023038224
$ws.Range("F9").Value = 0.3353220329595743
$ws.Range("G9").Value = 0.1906524392664082
$ws.Range("H9").Value = 0.3487996483136939
$ws.Range("M9").Value = 1.110542873411021
$ws.Range("O9").Value = 0.9851449114098614
# Row 10
$ws.Range("B10").Value = 0.9186792045018137
$ws.Range("C10").Value = 0.2416343682971842
$ws.Range("D10").Value = 0.03098469536905668
$ws.Range("F10").Value = 0.3428195408954551
$ws.Range("G10").Value = 0.1955364465317899
$ws.Range("H10").Value = 0.3452288411542241
$ws.Range("M10").Value = 1.296990622120717
$ws.Range("O10").Value = 0.9884139801868912
# Row 11
$ws.Range("B11").Value = 0.9822869227480169
$ws.Range("C11").Value = 0.2526820761409851
$ws.Range("D11").Value = 0.03313555504199428
$ws.Range("F11").Value = 0.346573929075717
$ws.Range("G11").Value = 0.1980237865991796
$ws.Range("H11").Value = 0.3438561581651953
$ws.Range("M11").Value = 1.382589326474104
$ws.Range("O11").Value = 0.9909820100158129
# Row 12
$ws.Range("B12").Value = 1.006332800113682
$ws.Range("C12").Value = 0.2568601848907974
$ws.Range("D12").Value = 0.03394877323182754
$ws.Range("F12").Value = 0.3480453193643527
$ws.Range("G12").Value = 0.1990042370663332
$ws.Range("H12").Value = 0.3433726362889189
$ws.Range("M12").Value = 1.415123620169084
$ws.Range("O12").Value = 0.99211085549328
# Row 13
$ws.Range("B13").Value = 1.001155944654556
$ws.Range("C13").Value = 0.2559606029435599
$ws.Range("D13").Value = 0.03377368940766701
$ws.Range("F13").Value = 0.3477262153137204
$ws.Range("G13").Value = 0.1987913593452362
$ws.Range("H13").Value = 0.3434751562300988
$ws.Range("M13").Value = 1.408111324865175
$ws.Range("O13").Value = 0.9918607678272338
# Row 14
$ws.Range("B14").Value = 0.9842660243407977
$ws.Range("C14").Value = 0.2530259229687886
$ws.Range("D14").Value = 0.03320248470774345
$ws.Range("F14").Value = 0.3466939840143155
$ws.Range("G14").Value = 0.1981036742225086
$ws.Range("H14").Value = 0.3438156508407531
$ws.Range("M14").Value = 1.385263496046122
$ws.Range("O14").Value = 0.9910717410966754
# Row 15
$ws.Range("B15").Value = 0.9739150639750846
$ws.Range("C15").Value = 0.2512276270403504
$ws.Range("D15").Value = 0.03285243869410692
$ws.Range("F15").Value = 0.3460681905810574
$ws.Range("G15").Value = 0.1976874780475555
$ws.Range("H15").Value = 0.3440289414588023
$ws.Range("M15").Value = 1.371284380741059
$ws.Range("O15").Value = 0.9906088341145676
# Row 16
$ws.Range("B16").Value = 0.9145166636460544
$ws.Range("C16").Value = 0.2409116422126658
$ws.Range("D16").Value = 0.03084395784217975
$ws.Range("F16").Value = 0.3425811215495784
$ws.Range("G16").Value = 0.1953792678380779
$ws.Range("H16").Value = 0.3453236216514455
$ws.Range("M16").Value = 1.291412920954059
$ws.Range("O16").Value = 0.9882679857481094
# Row 17
$ws.Range("B17").Value = 0.8780068000432948
$ws.Range("C17").Value = 0.2345739767123689
$ws.Range("D17").Value = 0.02960963088790436
$ws.Range("F17").Value = 0.3405301389628974
$ws.Range("G17").Value = 0.194031532562704
$ws.Range("H17").Value = 0.3461823951393654
$ws.Range("M17").Value = 1.242620017060005
$ws.Range("O17").Value = 0.9871094784401748
# Row 18
$ws.Range("B18").Value = 0.8569818898522499
$ws.Range("C18").Value = 0.2309255122476657
$ws.Range("D18").Value = 0.02889889386725741
$ws.Range("F18").Value = 0.3393828066621722
$ws.Range("G18").Value = 0.19328132720225
$ws.Range("H18").Value = 0.3467000242996079
$ws.Range("M18").Value = 1.214628846955435
$ws.Range("O18").Value = 0.986544819505653
# Row 19
$ws.Range("B19").Value = 0.8498589146711311
$ws.Range("C19").Value = 0.2296896686454488
$ws.Range("D19").Value = 0.02865811759949111
$ws.Range("F19").Value = 0.3389998860195931
$ws.Range("G19").Value = 0.1930315993263036
$ws.Range("H19").Value = 0.3468793498808154
$ws.Range("M19").Value = 1.205163871476302
$ws.Range("O19").Value = 0.9863710708814892
# Row 20
$ws.Range("B20").Value = 0.8818959808860427
$ws.Range("C20").Value = 0.2352489673844502
$ws.Range("D20").Value = 0.02974110876637326
$ws.Range("F20").Value = 0.3407451207749475
$ws.Range("G20").Value = 0.1941724137301009
$ws.Range("H20").Value = 0.3460885252395158
$ws.Range("M20").Value = 1.247806474436572
$ws.Range("O20").Value = 0.987222272717986
# Row 21
$ws.Range("B21").Value = 0.9892281315210312
$ws.Range("C21").Value = 0.2538880604469682
$ws.Range("D21").Value = 0.03337029621074805
$ws.Range("F21").Value = 0.3469958248940301
$ws.Range("G21").Value = 0.1983046149555037
$ws.Range("H21").Value = 0.3437146538826426
$ws.Range("M21").Value = 1.391971142831409
$ws.Range("O21").Value = 0.9912992454410414
# Row 22
$ws.Range("B22").Value = 1.059136086782303
$ws.Range("C22").Value = 0.2660379840068003
$ws.Range("D22").Value = 0.03573477433791084
$ws.Range("F22").Value = 0.3513707417286298
$ws.Range("G22").Value = 0.2012301050471592
$ws.Range("H22").Value = 0.3423747208106818
$ws.Range("M22").Value = 1.486893648217659
$ws.Range("O22").Value = 0.9948757756012867
# Row 23
$ws.Range("B23").Value = 1.021847494252995
$ws.Range("C23").Value = 0.2595564073192804
$ws.Range("D23").Value = 0.03447350564240992
$ws.Range("F23").Value = 0.3490091723600059
$ws.Range("G23").Value = 0.1996480230971258
$ws.Range("H23").Value = 0.343070483648404
$ws.Range("M23").Value = 1.436164948941226
$ws.Range("O23").Value = 0.9928831517464687
# Row 24
$ws.Range("B24").Value = 0.8801377907910819
$ws.Range("C24").Value = 0.2349438191915283
$ws.Range("D24").Value = 0.02968167105751007
$ws.Range("F24").Value = 0.3406478284519068
$ws.Range("G24").Value = 0.194108644678316
$ws.Range("H24").Value = 0.34613088935739
$ws.Range("M24").Value = 1.245461489104613
$ws.Range("O24").Value = 0.9871709627138614
# Row 25
$ws.Range("B25").Value = 0.7263174004750113
$ws.Range("C25").Value = 0.2082744169010198
$ws.Range("D25").Value = 0.02448312755733184
$ws.Range("F25").Value = 0.3328640302678352
$ws.Range("G25").Value = 0.1890869227953544
$ws.Range("H25").Value = 0.3503372018471893
$ws.Range("M25").Value = 1.042527182776013
$ws.Range("O25").Value = 0.9848930400205518

Write-Host "Updated 192 cells across rows 2-25 (380 kV case)"